$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bars")
$ws.Range("L1").EntireColumn.Insert()
Write-Host "done"
for ($c = 1; $c -le 17; $c++) {
    $r = $ws.Cells.Item(1, $c)
    $v = $r.Value()
    Write-Host "col $c = $v"
}
